$wb = $excel.ActiveWorkbook

# Row -> new value for column F (想去人数) on the "展览" and "全部类型" sheets
$updates = @{
    8  = 502
    9  = 6660
    10 = 181
    11 = 147
    12 = 1033
    13 = 370
    14 = 119
    15 = 184
    16 = 537
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
